# Weekly update: add two new rows of fresh data (week of 2022-05-25 / serial 44706)
# at the top of the "Vega Monumental Concepción - Lechuga" data block (old rows 633-653),
# pushing the existing rows down by two (they become rows 635-655).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above the current row 633.
$ws.Rows.Item(633).EntireRow.Insert()
$ws.Rows.Item(633).EntireRow.Insert()

# --- New row 633: Lechuga Marina ---
$ws.Cells.Item(633,1).Value  = 11
$ws.Cells.Item(633,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(633,3).Value  = "Bíobío"
$ws.Cells.Item(633,4).Value  = 44706
$ws.Cells.Item(633,5).Value  = 8
$ws.Cells.Item(633,6).Value  = 100112033
$ws.Cells.Item(633,7).Value  = "Lechuga"
$ws.Cells.Item(633,8).Value  = "Marina"
$ws.Cells.Item(633,9).Value  = "Primera"
$ws.Cells.Item(633,10).Value = 100
$ws.Cells.Item(633,11).Value = 6000
$ws.Cells.Item(633,12).Value = 6500
$ws.Cells.Item(633,13).Value = 6250
$ws.Cells.Item(633,14).Value = "`$/caja 15 unidades"
$ws.Cells.Item(633,15).Value = "Región Metropolitana"
$ws.Cells.Item(633,16).Value = 417
$ws.Cells.Item(633,17).Value = 15
$ws.Cells.Item(633,18).Value = "Hortaliza"

# --- New row 634: Lechuga Milanesa ---
$ws.Cells.Item(634,1).Value  = 11
$ws.Cells.Item(634,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(634,3).Value  = "Bíobío"
$ws.Cells.Item(634,4).Value  = 44706
$ws.Cells.Item(634,5).Value  = 8
$ws.Cells.Item(634,6).Value  = 100112033
$ws.Cells.Item(634,7).Value  = "Lechuga"
$ws.Cells.Item(634,8).Value  = "Milanesa"
$ws.Cells.Item(634,9).Value  = "Primera"
$ws.Cells.Item(634,10).Value = 100
$ws.Cells.Item(634,11).Value = 6000
$ws.Cells.Item(634,12).Value = 6500
$ws.Cells.Item(634,13).Value = 6250
$ws.Cells.Item(634,14).Value = "`$/caja 20 unidades"
$ws.Cells.Item(634,15).Value = "Región Metropolitana"
$ws.Cells.Item(634,16).Value = 312
$ws.Cells.Item(634,17).Value = 20
$ws.Cells.Item(634,18).Value = "Hortaliza"
